$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price/Volume columns so numeric-looking
# strings (e.g. "233.02", "0.0670") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.030.64"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.850.82"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "233.02"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "40.79"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").Value = "0.331"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "0.0986"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "2.119.67"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").Value = "1.850.17"
$ws.Range("D15").Value = "0.676"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "35.078.21"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "70.11"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").Value = "240.58"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "12.27"
$ws.Range("E21").Value = "  +4.42%  "
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +3.12%  "
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "173.02"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "7.87"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "17.53"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "0.0554"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").Value = "  +21.99%  "
$ws.Range("E35").Value = "  +11.78%  "
$ws.Range("D36").Value = "0.755"
$ws.Range("E36").Value = "  +10.97%  "
$ws.Range("E37").Value = "  +7.71%  "
$ws.Range("E38").Value = "  +13.16%  "
$ws.Range("D39").Value = "90.41"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.351.88"
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "14.65"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").Value = "2.037.77"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("E49").Value = "  +19.64%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "0.0670"
$ws.Range("E51").Value = "  -0.31%  "
